$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 1233.7
$ws.Range("J106").Value = 1696.6666
$ws.Range("L106").Value = 1696.6666
$ws.Range("N106").Value = -2958.6666
$ws.Range("H116").Value = 4289.7896
$ws.Range("I116").Value = 3490
$ws.Range("K116").Value = 3490
$ws.Range("M116").Value = -48
$ws.Range("H125").Value = 3166.6667
$ws.Range("I125").Value = 1700
$ws.Range("J125").Value = 5000
$ws.Range("K125").Value = 15300
$ws.Range("L125").Value = 45000
$ws.Range("M125").Value = -12840
$ws.Range("N125").Value = -49920

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2372.8333
$ws.Range("I63").Value = 2263.7778
$ws.Range("K63").Value = 2263.7778
$ws.Range("M63").Value = -1577.7778
$ws.Range("H66").Value = 2372.8333
$ws.Range("I66").Value = 2263.7778
$ws.Range("K66").Value = 11318.889
$ws.Range("M66").Value = -7886.888999999999
$ws.Range("H88").Value = 2088.7273
$ws.Range("I88").Value = 1747
$ws.Range("J88").Value = 3000
$ws.Range("K88").Value = 1747
$ws.Range("L88").Value = 3000
$ws.Range("M88").Value = -1341
$ws.Range("N88").Value = -3812
$ws.Range("H91").Value = 2088.7273
$ws.Range("I91").Value = 1747
$ws.Range("J91").Value = 3000
$ws.Range("K91").Value = 1747
$ws.Range("L91").Value = 3000
$ws.Range("M91").Value = -343
$ws.Range("N91").Value = -5808
$ws.Range("H103").Value = 25098.275
$ws.Range("J103").Value = 25098.275
$ws.Range("L103").Value = 25098.275
$ws.Range("N103").Value = -27442.275
$ws.Range("H107").Value = 50000
$ws.Range("J107").Value = 50000
$ws.Range("L107").Value = 50000
$ws.Range("N107").Value = -57680

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 44633.918
$ws.Range("I86").Value = 1257.1428
$ws.Range("J86").Value = 105361.4
$ws.Range("K86").Value = 1257.1428
$ws.Range("L86").Value = 105361.4
$ws.Range("M86").Value = -134.1428000000001
$ws.Range("N86").Value = -107607.4
$ws.Range("H89").Value = 44633.918
$ws.Range("I89").Value = 1257.1428
$ws.Range("J89").Value = 105361.4
$ws.Range("K89").Value = 6285.714
$ws.Range("L89").Value = 526807
$ws.Range("M89").Value = -669.7139999999999
$ws.Range("N89").Value = -538039
$ws.Range("H99").Value = 3398.5715
$ws.Range("I99").Value = 1945
$ws.Range("K99").Value = 1945
$ws.Range("M99").Value = -447

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2132010
$ws.Range("I31").Value = 3128459
$ws.Range("J31").Value = 6252.2666
$ws.Range("K31").Value = 3128459
$ws.Range("L31").Value = 6252.2666
$ws.Range("M31").Value = -3128164
$ws.Range("N31").Value = -6842.2666
$ws.Range("H34").Value = 2132010
$ws.Range("I34").Value = 3128459
$ws.Range("J34").Value = 6252.2666
$ws.Range("K34").Value = 3128459
$ws.Range("L34").Value = 6252.2666
$ws.Range("M34").Value = -3128257
$ws.Range("N34").Value = -6656.2666
$ws.Range("H58").Value = 8622881
$ws.Range("I58").Value = 1416.6154
$ws.Range("J58").Value = 26319570
$ws.Range("K58").Value = 1416.6154
$ws.Range("L58").Value = 26319570
$ws.Range("M58").Value = -1213.6154
$ws.Range("N58").Value = -26319976
$ws.Range("H62").Value = 3616.875
$ws.Range("I62").Value = 2781.4285
$ws.Range("J62").Value = 4266.6665
$ws.Range("K62").Value = 2781.4285
$ws.Range("L62").Value = 4266.6665
$ws.Range("M62").Value = -2157.4285
$ws.Range("N62").Value = -5514.6665
$ws.Range("H65").Value = 3616.875
$ws.Range("I65").Value = 2781.4285
$ws.Range("J65").Value = 4266.6665
$ws.Range("K65").Value = 13907.1425
$ws.Range("L65").Value = 21333.3325
$ws.Range("M65").Value = -10787.1425
$ws.Range("N65").Value = -27573.3325
$ws.Range("H107").Value = 1778.5
$ws.Range("I107").Value = 279.27274
$ws.Range("J107").Value = 3610.889
$ws.Range("K107").Value = 279.27274
$ws.Range("L107").Value = 3610.889
$ws.Range("M107").Value = 1640.72726
$ws.Range("N107").Value = -7450.889
$ws.Range("H132").Value = 2170.2886
$ws.Range("I132").Value = 1567.4166
$ws.Range("K132").Value = 4702.2498
$ws.Range("M132").Value = -2172.2498
$ws.Range("H134").Value = 1218.0878
$ws.Range("I134").Value = 866.43335
$ws.Range("J134").Value = 1608.8148
$ws.Range("K134").Value = 2599.30005
$ws.Range("L134").Value = 4826.4444
$ws.Range("M134").Value = -64.30004999999983
$ws.Range("N134").Value = -9896.4444
$ws.Range("H136").Value = 8622881
$ws.Range("I136").Value = 1416.6154
$ws.Range("J136").Value = 26319570
$ws.Range("K136").Value = 4249.8462
$ws.Range("L136").Value = 78958710
$ws.Range("M136").Value = -1699.8462
$ws.Range("N136").Value = -78963810

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 11000
$ws.Range("I87").Value = 9000
$ws.Range("J87").Value = 15000
$ws.Range("K87").Value = 27000
$ws.Range("L87").Value = 45000
$ws.Range("M87").Value = -25752
$ws.Range("N87").Value = -47496
$ws.Range("H90").Value = 11000
$ws.Range("I90").Value = 9000
$ws.Range("J90").Value = 15000
$ws.Range("K90").Value = 81000
$ws.Range("L90").Value = 135000
$ws.Range("M90").Value = -74760
$ws.Range("N90").Value = -147480
$ws.Range("H107").Value = 1869.8
$ws.Range("J107").Value = 3899.3333
$ws.Range("L107").Value = 11697.9999
$ws.Range("N107").Value = -15537.9999
$ws.Range("H132").Value = 960
$ws.Range("I132").Value = 960
$ws.Range("K132").Value = 8640
$ws.Range("M132").Value = -6110

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4320.222
$ws.Range("I70").Value = 4605.4546
$ws.Range("J70").Value = 3872
$ws.Range("K70").Value = 4605.4546
$ws.Range("L70").Value = 3872
$ws.Range("M70").Value = -4335.4546
$ws.Range("N70").Value = -4412
$ws.Range("H73").Value = 4320.222
$ws.Range("I73").Value = 4605.4546
$ws.Range("J73").Value = 3872
$ws.Range("K73").Value = 4605.4546
$ws.Range("L73").Value = 3872
$ws.Range("M73").Value = -3669.4546
$ws.Range("N73").Value = -5744
$ws.Range("H126").Value = 2522.9688
$ws.Range("I126").Value = 1438.75
$ws.Range("J126").Value = 3607.1875
$ws.Range("K126").Value = 4316.25
$ws.Range("L126").Value = 10821.5625
$ws.Range("M126").Value = -1846.25
$ws.Range("N126").Value = -15761.5625
